# Apply Amharic translations to "Facilitator Guideline - Airport Problem.docx"
#
# Each English string that was translated in the source commit is unique
# inside the document body except for "Geometry", which also appears a
# second time inside an unrelated "VIDEO PAUSE\nGeometry" run further
# down in the document. Using wdReplaceOne (replace only the first match)
# together with a search that always starts at the top of the document
# ensures we only touch the "Geometry" that is the value for the "Topic"
# row, leaving the other "Geometry" occurrence untouched - exactly as the
# diff requires.

$d = $word.ActiveDocument

function Replace-FirstMatch($findText, $replacementText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute(
        $findText,      # FindText
        $true,          # MatchCase
        $true,          # MatchWholeWord
        $false,         # MatchWildcards
        $false,         # MatchSoundsLike
        $false,         # MatchAllWordForms
        $true,          # Forward
        1,              # Wrap (wdFindContinue)
        $false,         # Format
        $replacementText, # ReplaceWith
        1               # Replace (wdReplaceOne -> only the first hit)
    ) | Out-Null
}

Replace-FirstMatch "Video Title" "የቪድዮ አርዕስት"
Replace-FirstMatch "The Airport Problem" "የአየር መንገድ ጥያቄ"
Replace-FirstMatch "Topic" "የትምህርት ርዕስ"
Replace-FirstMatch "Geometry" "ጂኦሜትሪ"
Replace-FirstMatch "Aim(s)" "አላማ (ዎች)"
Replace-FirstMatch "Get the intuitive idea of a minimization problem, figure out how to practically implement minimization problems." "የወጭ ቅነሳን መነሣ ሃሳብ ለማገኘት፣ የወጭ ቅነሳ ጥያቄዎችን እንደት በተግባር መተግበር አነደሚቻል ለማየት"
Replace-FirstMatch "Length" "ርዝመት"
